$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 637.41174
$ws.Range("J33").Value = 65
$ws.Range("L33").Value = 65
$ws.Range("N33").Value = -523
$ws.Range("H64").Value = 6192.857
$ws.Range("I64").Value = 5154.5454
$ws.Range("J64").Value = 10000
$ws.Range("K64").Value = 5154.5454
$ws.Range("L64").Value = 10000
$ws.Range("M64").Value = -4906.5454
$ws.Range("N64").Value = -10496
$ws.Range("H67").Value = 6192.857
$ws.Range("I67").Value = 5154.5454
$ws.Range("J67").Value = 10000
$ws.Range("K67").Value = 5154.5454
$ws.Range("L67").Value = 10000
$ws.Range("M67").Value = -4296.5454
$ws.Range("N67").Value = -11716
$ws.Range("H116").Value = 5580.4443
$ws.Range("J116").Value = 5761.2856
$ws.Range("L116").Value = 5761.2856
$ws.Range("N116").Value = -12645.2856
$ws.Range("H132").Value = 15153639
$ws.Range("I132").Value = 17243580
$ws.Range("J132").Value = 1559.875
$ws.Range("K132").Value = 51730740
$ws.Range("L132").Value = 4679.625
$ws.Range("M132").Value = -51728210
$ws.Range("N132").Value = -9739.625
$ws.Range("H135").Value = 1793.8788
$ws.Range("I135").Value = 662.1
$ws.Range("K135").Value = 5958.900000000001
$ws.Range("M135").Value = -3423.900000000001
$ws.Range("H137").Value = 3061.1853
$ws.Range("I137").Value = 2754.6924
$ws.Range("J137").Value = 3345.7856
$ws.Range("K137").Value = 8264.0772
$ws.Range("L137").Value = 10037.3568
$ws.Range("M137").Value = -5714.0772
$ws.Range("N137").Value = -15137.3568
$ws.Range("H138").Value = 2002.1855
$ws.Range("I138").Value = 957.3415
$ws.Range("J138").Value = 2767.1606
$ws.Range("K138").Value = 2872.0245
$ws.Range("L138").Value = 8301.481800000001
$ws.Range("M138").Value = 2267.9755
$ws.Range("N138").Value = -18581.4818
$ws.Range("H141").Value = 1408.3143
$ws.Range("J141").Value = 2414.75
$ws.Range("L141").Value = 7244.25
$ws.Range("N141").Value = -17604.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3086.5269
$ws.Range("I32").Value = 2258.96
$ws.Range("K32").Value = 2258.96
$ws.Range("M32").Value = -1971.96
$ws.Range("H61").Value = 1361.5135
$ws.Range("I61").Value = 1260.5834
$ws.Range("K61").Value = 1260.5834
$ws.Range("M61").Value = -1048.5834
$ws.Range("H74").Value = 254673.55
$ws.Range("I74").Value = 131563.62
$ws.Range("K74").Value = 131563.62
$ws.Range("M74").Value = -130689.62
$ws.Range("H77").Value = 254673.55
$ws.Range("I77").Value = 131563.62
$ws.Range("K77").Value = 657818.1
$ws.Range("M77").Value = -653450.1
$ws.Range("H110").Value = 1389728.6
$ws.Range("I110").Value = 1852751
$ws.Range("J110").Value = 661.6
$ws.Range("K110").Value = 1852751
$ws.Range("L110").Value = 661.6
$ws.Range("M110").Value = -1850706
$ws.Range("N110").Value = -4751.6
$ws.Range("H122").Value = 634034.8
$ws.Range("I122").Value = 2210.4167
$ws.Range("J122").Value = 2318900
$ws.Range("K122").Value = 6631.250100000001
$ws.Range("L122").Value = 6956700
$ws.Range("M122").Value = -4181.250100000001
$ws.Range("N122").Value = -6961600
$ws.Range("H132").Value = 1362.0588
$ws.Range("I132").Value = 908.91113
$ws.Range("J132").Value = 4760.6665
$ws.Range("K132").Value = 2726.73339
$ws.Range("L132").Value = 14281.9995
$ws.Range("M132").Value = -196.7333899999999
$ws.Range("N132").Value = -19341.9995
$ws.Range("H136").Value = 1361.5135
$ws.Range("I136").Value = 1260.5834
$ws.Range("K136").Value = 3781.7502
$ws.Range("M136").Value = -1231.7502

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 17985558
$ws.Range("I99").Value = 47953710
$ws.Range("J99").Value = 4665.8
$ws.Range("K99").Value = 47953710
$ws.Range("L99").Value = 4665.8
$ws.Range("M99").Value = -47952212
$ws.Range("N99").Value = -7661.8
$ws.Range("H107").Value = 4763757.5
$ws.Range("I107").Value = 5496413
$ws.Range("J107").Value = 1495.5
$ws.Range("K107").Value = 5496413
$ws.Range("L107").Value = 1495.5
$ws.Range("M107").Value = -5494493
$ws.Range("N107").Value = -5335.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2303.423
$ws.Range("I31").Value = 1505.05
$ws.Range("J31").Value = 4964.6665
$ws.Range("K31").Value = 1505.05
$ws.Range("L31").Value = 4964.6665
$ws.Range("M31").Value = -1210.05
$ws.Range("N31").Value = -5554.6665
$ws.Range("H34").Value = 2303.423
$ws.Range("I34").Value = 1505.05
$ws.Range("J34").Value = 4964.6665
$ws.Range("K34").Value = 1505.05
$ws.Range("L34").Value = 4964.6665
$ws.Range("M34").Value = -1303.05
$ws.Range("N34").Value = -5368.6665
$ws.Range("H58").Value = 5144.4736
$ws.Range("I58").Value = 5283.6665
$ws.Range("J58").Value = 4622.5
$ws.Range("K58").Value = 5283.6665
$ws.Range("L58").Value = 4622.5
$ws.Range("M58").Value = -5080.6665
$ws.Range("N58").Value = -5028.5
$ws.Range("H122").Value = 2666.2632
$ws.Range("I122").Value = 2405.3333
$ws.Range("J122").Value = 3113.5715
$ws.Range("K122").Value = 7215.999899999999
$ws.Range("L122").Value = 9340.7145
$ws.Range("M122").Value = -4765.999899999999
$ws.Range("N122").Value = -14240.7145
$ws.Range("H132").Value = 2499.182
$ws.Range("I132").Value = 1719.6666
$ws.Range("K132").Value = 5158.9998
$ws.Range("M132").Value = -2628.9998
$ws.Range("H136").Value = 5144.4736
$ws.Range("I136").Value = 5283.6665
$ws.Range("J136").Value = 4622.5
$ws.Range("K136").Value = 15850.9995
$ws.Range("L136").Value = 13867.5
$ws.Range("M136").Value = -13300.9995
$ws.Range("N136").Value = -18967.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 417734.88
$ws.Range("J46").Value = 1293.6666
$ws.Range("L46").Value = 3880.9998
$ws.Range("N46").Value = -4062.9998
$ws.Range("H125").Value = 833.3333
$ws.Range("I125").Value = 500
$ws.Range("K125").Value = 1500
$ws.Range("M125").Value = 3420

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4098.533
$ws.Range("J122").Value = 6213.6665
$ws.Range("L122").Value = 18640.9995
$ws.Range("N122").Value = -23540.9995
$ws.Range("H132").Value = 2693.3953
$ws.Range("I132").Value = 2141.0881
$ws.Range("J132").Value = 4779.8887
$ws.Range("K132").Value = 6423.2643
$ws.Range("L132").Value = 14339.6661
$ws.Range("M132").Value = -3893.2643
$ws.Range("N132").Value = -19399.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 840.0769
$ws.Range("J22").Value = 942.8570999999999
$ws.Range("L22").Value = 942.8570999999999
$ws.Range("N22").Value = -1532.8571
$ws.Range("H27").Value = 840.0769
$ws.Range("J27").Value = 942.8570999999999
$ws.Range("L27").Value = 942.8570999999999
$ws.Range("N27").Value = -1156.8571
$ws.Range("H122").Value = 7404.4443
$ws.Range("I122").Value = 4614.5
$ws.Range("J122").Value = 9636.4
$ws.Range("K122").Value = 13843.5
$ws.Range("L122").Value = 28909.2
$ws.Range("M122").Value = -11393.5
$ws.Range("N122").Value = -33809.2
$ws.Range("H132").Value = 7068.3335
$ws.Range("I132").Value = 7076.2905
$ws.Range("K132").Value = 21228.8715
$ws.Range("M132").Value = -18698.8715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 62507490
$ws.Range("I107").Value = 90913800
$ws.Range("J107").Value = 13600.2
$ws.Range("K107").Value = 272741400
$ws.Range("L107").Value = 40800.60000000001
$ws.Range("M107").Value = -272739480
$ws.Range("N107").Value = -44640.60000000001
